$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text edits -------------------------------------------------
# D22 keeps referencing the same shared-string slot; only its text changes.
$ws.Range("D22").Value = "19, 21, 25, 27, 28, 29, 52, 80"

# D23 stops being a text cell and becomes numeric (page number 13). This frees
# up the shared-string slot it used to hold ("19, 21, 25, 27, 29, ").
$ws.Range("D23").Value = 13

# G22 is a brand-new note cell; it ends up reusing the slot D23 just freed.
$ws.Range("G22").Value = "A lo largo del documento se cambio 'sistema' por 'software' donde se necesitaba"

# --- New date values -----------------------------------------------------------
$ws.Range("E22").Value = 43269
$ws.Range("E22").NumberFormat = "d-mmm"

$ws.Range("E23").Value = 43269
$ws.Range("E23").NumberFormat = "d-mmm"

# --- Re-use F2's existing fill style on F23 (formats only, keep its value) ----
$ws.Range("F2").Copy()
$ws.Range("F23").PasteSpecial(-4122)

# --- Row 22 grows taller to fit the new "Hecho" note --------------------------
$ws.Rows.Item(22).RowHeight = 57

# --- Move the active selection to F22 ------------------------------------------
$ws.Range("F22").Select()
